$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 91813
$ws.Range("B4").Value = 92467
$ws.Range("B5").Value = 91870
$ws.Range("B6").Value = 92301
